# Add a new weekly-report row (row 17) to the "下周工作" table on the
# right-hand side of the sheet, mirroring rows 13-16, and move the
# active-cell selection past the newly inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G17").Value = "微信客户端界面Mock图"
$ws.Range("H17").Value = 20160710
$ws.Range("I17").Value = "魏鑫"

$ws.Range("G18").Select()
